$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns keep their text representation (values like
# "588.03" or "181.48" would otherwise be auto-coerced to numbers by Excel).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "69.348.54"
$ws.Range("D3").Value = "3.407.70"
$ws.Range("E3").Value = "  +2.65%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "588.03"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").Value = "181.48"
$ws.Range("E6").Value = "  +4.33%  "
$ws.Range("E7").Value = "  +2.14%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").Value = "0.199"
$ws.Range("E9").Value = "  +10.21%  "
$ws.Range("E10").Value = "  +3.04%  "
$ws.Range("D11").Value = "48.63"
$ws.Range("E11").Value = "  +4.31%  "
$ws.Range("E12").Value = "  +5.11%  "
$ws.Range("D13").Value = "686.06"
$ws.Range("E13").Value = "  -1.57%  "
$ws.Range("D14").Value = "8.70"
$ws.Range("E14").Value = "  +4.27%  "
$ws.Range("D15").Value = "3.957.04"
$ws.Range("E15").Value = "  +2.58%  "
$ws.Range("D16").Value = "69.467.33"
$ws.Range("E16").Value = "  +2.80%  "
$ws.Range("D17").Value = "3.409.00"
$ws.Range("E18").Value = "  +1.75%  "
$ws.Range("D19").Value = "17.77"
$ws.Range("E19").Value = "  +1.99%  "
$ws.Range("D20").Value = "11.34"
$ws.Range("E20").Value = "  +2.42%  "
$ws.Range("E21").Value = "  +2.82%  "
$ws.Range("D22").Value = "17.25"
$ws.Range("E22").Value = "  +2.68%  "
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("D24").Value = "103.52"
$ws.Range("E24").Value = "  +1.99%  "
$ws.Range("D25").Value = "3.94"
$ws.Range("E25").Value = "  +1.47%  "
$ws.Range("D26").Value = "2.73"
$ws.Range("E26").Value = "  +2.51%  "
$ws.Range("D27").Value = "9.72"
$ws.Range("E27").Value = "  +4.26%  "
$ws.Range("D29").Value = "8.82"
$ws.Range("E29").Value = "  +4.14%  "
$ws.Range("D30").Value = "6.97"
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("D31").Value = "563.45"
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("D32").Value = "11.17"
$ws.Range("E32").Value = "  +1.99%  "
$ws.Range("B33").Value = "dogwifhat"
$ws.Range("C33").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D33").Value = "3.61"
$ws.Range("E33").Value = "  +11.78%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "0.107"
$ws.Range("E34").Value = "  +1.91%  "
$ws.Range("D35").Value = "58.63"
$ws.Range("E35").Value = "  +3.76%  "
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").Value = "3.666.61"
$ws.Range("E37").Value = "  -0.99%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.142"
$ws.Range("E38").Value = "  +7.58%  "
$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D39").Value = "36.07"
$ws.Range("E39").Value = "  +3.57%  "
$ws.Range("D40").Value = "0.0₃0724"
$ws.Range("E40").Value = "  +9.14%  "
$ws.Range("E41").Value = "  +5.02%  "
$ws.Range("E42").Value = "  +3.61%  "
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").Value = "0.340"
$ws.Range("E43").Value = "  +2.30%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0427"
$ws.Range("E44").Value = "  +5.66%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "3.34"
$ws.Range("E45").Value = "  +1.29%  "
$ws.Range("E46").Value = "  +2.20%  "
$ws.Range("E47").Value = "  +1.60%  "
$ws.Range("E48").Value = "  +6.00%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").Value = "132.68"
$ws.Range("E50").Value = "  +1.33%  "
$ws.Range("D51").Value = "2.66"
$ws.Range("E51").Value = "  +2.86%  "
